$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "Fabio"
$ws.Range("A5").Value = "Ellen"
$ws.Range("A6").Value = "Fabio"
$ws.Range("A7").Value = "Paritosh"
$ws.Range("A8").Value = "Tina"
$ws.Range("A9").Value = "Paritosh"
$ws.Range("A10").Value = "Ellen"
$ws.Range("A11").Value = "Tina"

$ws.Range("G6").Select()
